$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Noah"
$ws.Range("C2").Value = "Seligson"
$ws.Range("G2").Value = "Education, Environment, Sports & Recreation, Coding & Software Development, Music & Performance, Health & Wellness, Animal Welfare"
